# Atualizando para visitar clientes - by matheus
#
# Insere 4 novos dias (15-18) no mes de maio/2025 (linhas 16-19), empurrando
# todos os registros subsequentes (meses 04, 03 e 02 de 2025) quatro linhas
# para baixo, e atualiza os valores de faturamento dos dias 13 e 14 de maio.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Atualiza os totais de faturamento ja existentes para os dias 13 e 14 de maio/2025
$ws.Range("B14").Value = 26531.86
$ws.Range("B15").Value = 36574.18

# Insere 4 linhas em branco antes da linha 16, empurrando o restante dos dados
# (meses 04, 03 e 02/2025) para baixo
$ws.Rows("16:19").Insert()

# Preenche as novas linhas com os dias 15 a 18 de maio/2025
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 33940.79
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2025
$ws.Range("E16").Value = "05/2025"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 30403.76
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2025
$ws.Range("E17").Value = "05/2025"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 14533.8
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = "05/2025"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 8085.01
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = "05/2025"
